# Appended data to sheet 'Sheet1'
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A holds a date-like string ("2025-10-17"). Format the cell as
# Text first so Excel stores the literal string instead of silently
# converting it to a date serial number, then restore the default
# "Normal" style so the cell keeps the workbook's plain/default
# formatting (matching the other plain-text cells on the sheet).
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "2025-10-17"
$ws.Cells.Item(3, 1).Style = "Normal"

$ws.Cells.Item(3, 2).Value = "ZZZ"
$ws.Cells.Item(3, 3).Value = "456CDX009"
$ws.Cells.Item(3, 4).Value = "Anna Nagar"
